$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column L (header + data + blank rows) over to the
# new column M so it inherits the same borders/fills/fonts/number formats.
$ws.Range("L1:L10").Copy()
$ws.Range("M1:M10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new "transaction_status_blockchain" column.
$ws.Range("M1").Value = "transaction_status_blockchain"
$ws.Range("M2").Value = "transaction_status_blockchain_link"
$ws.Range("M3").Value = "www_link"

# Widen column M to fit the new header/content (target stored width 42.3906;
# the engine quantizes ColumnWidth to 1/6-character steps before adding the
# fixed 5/6 padding baked into the stored OOXML `width`, so 41.5 is the
# closest achievable input).
$ws.Columns.Item(13).ColumnWidth = 41.5
